$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I4").Value = 0.2249675675517139
$ws.Range("J4").Value = 0.5161496630566509
$ws.Range("K4").Value = 0.765648256423696
$ws.Range("L4").Value = 3.108089542562038
